$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 308
$ws.Range("F4").Value = 420
$ws.Range("F5").Value = 8540
$ws.Range("F6").Value = 10
$ws.Range("F7").Value = 10691
$ws.Range("F13").Value = 116
$ws.Range("F20").Value = 411
$ws.Range("F22").Value = 1816
$ws.Range("F23").Value = 78
$ws.Range("F24").Value = 552
$ws.Range("F25").Value = 343
$ws.Range("F27").Value = 63
$ws.Range("F29").Value = 58
$ws.Range("F30").Value = 1183
$ws.Range("F33").Value = 1414
$ws.Range("F34").Value = 444
$ws.Range("F36").Value = 285
$ws.Range("F37").Value = 21
$ws.Range("F38").Value = 129
$ws.Range("F39").Value = 510
$ws.Range("F40").Value = 347
$ws.Range("F41").Value = 94
$ws.Range("F42").Value = 789
$ws.Range("F43").Value = 642
$ws.Range("F45").Value = 103
$ws.Range("F46").Value = 99

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 44
$ws.Range("F9").Value = 11
$ws.Range("F16").Value = 48

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2802

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 308
$ws.Range("F9").Value = 420
$ws.Range("F10").Value = 8540
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 10691
$ws.Range("F15").Value = 116
$ws.Range("F19").Value = 1816
$ws.Range("F20").Value = 78
$ws.Range("F21").Value = 552
$ws.Range("F23").Value = 63
$ws.Range("F26").Value = 44
$ws.Range("F28").Value = 1183
$ws.Range("F30").Value = 11
$ws.Range("F34").Value = 1414
$ws.Range("F35").Value = 444
$ws.Range("F38").Value = 129
$ws.Range("F39").Value = 510
$ws.Range("F41").Value = 347
$ws.Range("F42").Value = 94
$ws.Range("F43").Value = 789
$ws.Range("F45").Value = 48
$ws.Range("F47").Value = 642
$ws.Range("F48").Value = 103
$ws.Range("F49").Value = 99
